# Day 02 workbook touch-up:
#  - rename the tab from the generic "Sheet1" to the dated "18 Jan 2025"
#  - move the visible selection to B7 (and scroll the window toward row 10)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet tab
$ws.Name = "18 Jan 2025"

# Make sure this sheet/window is the active one before touching selection/scroll
$ws.Activate()

# Update the visible scroll position (top-left visible cell -> around A10)
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1

# Move the active selection to B7
$ws.Range("B7").Select()
